$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1:E36").Copy()
$ws.Range("B1").PasteSpecial()
